$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price/Volume columns hold plain text values in the source data
# (e.g. "278.28", "1.08%") rather than real numbers/percentages, so the
# affected cells are formatted as Text before the new values are written.
# This avoids Excel auto-converting them into numeric/percentage values.

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "278.28"
$ws.Range("E2").Value = "1.08%"
$rng.Style = "Normal"

$rng = $ws.Range("E3")
$rng.NumberFormat = "@"
$rng.Value = "2.54%"
$rng.Style = "Normal"

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "4.875"
$ws.Range("E4").Value = "-0.14%"
$rng.Style = "Normal"

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.06431"
$ws.Range("E5").Value = "1.37%"
$rng.Style = "Normal"

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "7.003"
$ws.Range("E6").Value = "1.35%"
$rng.Style = "Normal"

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "1.190"
$ws.Range("E7").Value = "-6.10%"
$rng.Style = "Normal"

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "0.8858"
$ws.Range("E8").Value = "1.89%"
$rng.Style = "Normal"

$rng = $ws.Range("E9")
$rng.NumberFormat = "@"
$rng.Value = "1.47%"
$rng.Style = "Normal"

$rng = $ws.Range("E10")
$rng.NumberFormat = "@"
$rng.Value = "1.77%"
$rng.Style = "Normal"

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.07497"
$ws.Range("E11").Value = "1.00%"
$rng.Style = "Normal"

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.02889"
$ws.Range("E12").Value = "-2.32%"
$rng.Style = "Normal"

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "0.08968"
$ws.Range("E13").Value = "-0.85%"
$rng.Style = "Normal"

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.001562"
$ws.Range("E14").Value = "-0.60%"
$rng.Style = "Normal"

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "0.0006348"
$ws.Range("E15").Value = "0.41%"
$rng.Style = "Normal"

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = "0.006159"
$ws.Range("E16").Value = "2.35%"
$rng.Style = "Normal"

$rng = $ws.Range("D17:E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = "3.480"
$ws.Range("E17").Value = "0.98%"
$rng.Style = "Normal"

$rng = $ws.Range("D18:E18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = "3.313"
$ws.Range("E18").Value = "-0.01%"
$rng.Style = "Normal"

$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$rng.Value = "0.08%"
$rng.Style = "Normal"

$rng = $ws.Range("E21")
$rng.NumberFormat = "@"
$rng.Value = "1.88%"
$rng.Style = "Normal"

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = "3.898"
$ws.Range("E22").Value = "-0.25%"
$rng.Style = "Normal"

$rng = $ws.Range("E23")
$rng.NumberFormat = "@"
$rng.Value = "0.84%"
$rng.Style = "Normal"

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = "0.001176"
$ws.Range("E25").Value = "-0.20%"
$rng.Style = "Normal"

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = "0.003879"
$ws.Range("E26").Value = "-7.93%"
$rng.Style = "Normal"

$rng = $ws.Range("E28")
$rng.NumberFormat = "@"
$rng.Value = "-1.68%"
$rng.Style = "Normal"

$rng = $ws.Range("E29")
$rng.NumberFormat = "@"
$rng.Value = "-1.72%"
$rng.Style = "Normal"

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.04151"
$ws.Range("E40").Value = "1.11%"
$rng.Style = "Normal"

$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = "0.006809"
$ws.Range("E41").Value = "-2.59%"
$rng.Style = "Normal"

$rng = $ws.Range("E42")
$rng.NumberFormat = "@"
$rng.Value = "0.34%"
$rng.Style = "Normal"

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.01174"
$ws.Range("E44").Value = "8.86%"
$rng.Style = "Normal"

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.00005311"
$ws.Range("E45").Value = "0.31%"
$rng.Style = "Normal"

$rng = $ws.Range("D46:E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = "1.687"
$ws.Range("E46").Value = "13.51%"
$rng.Style = "Normal"

$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$rng.Value = "-11.80%"
$rng.Style = "Normal"
